$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Baseline Training
$ws.Cells.Item(11, 2).Value = 45993
$ws.Cells.Item(11, 3).Value = 0.625
$ws.Cells.Item(11, 4).Value = 1
$ws.Range("E11").Formula = "=D11-C11"
$ws.Cells.Item(11, 6).Value = "Baseline Training"

# Row 12: Naive Quantization
$ws.Cells.Item(12, 2).Value = 45995
$ws.Cells.Item(12, 3).Value = 0.46875
$ws.Cells.Item(12, 4).Value = 0.53125
$ws.Range("E12").Formula = "=D12-C12"
$ws.Cells.Item(12, 6).Value = "Naive Quantization"

# Row 13: Naive Quantization
$ws.Cells.Item(13, 2).Value = 45995
$ws.Cells.Item(13, 3).Value = 0.59375
$ws.Cells.Item(13, 4).Value = 0.60416666666666663
$ws.Range("E13").Formula = "=D13-C13"
$ws.Cells.Item(13, 6).Value = "Naive Quantization"

# Apply the same time number format used by existing E column cells (E10)
# so the new cells reuse the existing style (matches the escaped format
# code stored in styles.xml) instead of minting a near-duplicate numFmt.
$ws.Range("E11:E13").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"

# Move the selection to mirror the author's final cursor position
$ws.Range("D13").Select()
